$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4825.6
$ws.Range("J62").Value = 5469
$ws.Range("L62").Value = 5469
$ws.Range("N62").Value = -6717
$ws.Range("H64").Value = 41670320
$ws.Range("J64").Value = 250003820
$ws.Range("L64").Value = 250003820
$ws.Range("N64").Value = -250004316
$ws.Range("H65").Value = 4825.6
$ws.Range("J65").Value = 5469
$ws.Range("L65").Value = 27345
$ws.Range("N65").Value = -33585
$ws.Range("H67").Value = 41670320
$ws.Range("J67").Value = 250003820
$ws.Range("L67").Value = 250003820
$ws.Range("N67").Value = -250005536
$ws.Range("H74").Value = 10031.654
$ws.Range("I74").Value = 11838.053
$ws.Range("K74").Value = 11838.053
$ws.Range("M74").Value = -10902.053
$ws.Range("H77").Value = 10031.654
$ws.Range("I77").Value = 11838.053
$ws.Range("K77").Value = 59190.265
$ws.Range("M77").Value = -54510.265
$ws.Range("H135").Value = 1594.8695
$ws.Range("J135").Value = 1819
$ws.Range("L135").Value = 16371
$ws.Range("N135").Value = -21441
$ws.Range("H138").Value = 4078.6072
$ws.Range("I138").Value = 5616.8184
$ws.Range("K138").Value = 16850.4552
$ws.Range("M138").Value = -11710.4552

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6831.3477
$ws.Range("I32").Value = 6164.316
$ws.Range("K32").Value = 6164.316
$ws.Range("M32").Value = -5877.316
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
$ws.Range("H63").Value = 1207.8572
$ws.Range("I63").Value = 1207.8572
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1207.8572
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -521.8571999999999
$ws.Range("N63").Value = $null
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H66").Value = 1207.8572
$ws.Range("I66").Value = 1207.8572
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 6039.286
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -2607.286
$ws.Range("N66").Value = -12950.4284
$ws.Range("H123").Value = 50413
$ws.Range("J123").Value = 50413
$ws.Range("L123").Value = 50413
$ws.Range("N123").Value = -60213

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17248130
$ws.Range("I20").Value = 23817172
$ws.Range("J20").Value = 4394
$ws.Range("K20").Value = 23817172
$ws.Range("L20").Value = 4394
$ws.Range("M20").Value = -23816925
$ws.Range("N20").Value = -4888
$ws.Range("H86").Value = 2283.2632
$ws.Range("I86").Value = 1753.5385
$ws.Range("K86").Value = 1753.5385
$ws.Range("M86").Value = -630.5385000000001
$ws.Range("H89").Value = 2283.2632
$ws.Range("I89").Value = 1753.5385
$ws.Range("K89").Value = 8767.692500000001
$ws.Range("M89").Value = -3151.692500000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6069
$ws.Range("I7").Value = 6785.125
$ws.Range("K7").Value = 6785.125
$ws.Range("M7").Value = -6672.125
$ws.Range("H22").Value = 1556.0834
$ws.Range("I22").Value = 1784.5
$ws.Range("J22").Value = 1099.25
$ws.Range("K22").Value = 1784.5
$ws.Range("L22").Value = 1099.25
$ws.Range("M22").Value = -1434.5
$ws.Range("N22").Value = -1799.25
$ws.Range("H99").Value = 9132.583000000001
$ws.Range("I99").Value = 10888
$ws.Range("J99").Value = 3866.3333
$ws.Range("K99").Value = 10888
$ws.Range("L99").Value = 3866.3333
$ws.Range("M99").Value = -9390
$ws.Range("N99").Value = -6862.3333
$ws.Range("H107").Value = 851.2143
$ws.Range("I107").Value = 586.5
$ws.Range("K107").Value = 586.5
$ws.Range("M107").Value = 1333.5
$ws.Range("H126").Value = 9132.583000000001
$ws.Range("I126").Value = 10888
$ws.Range("J126").Value = 3866.3333
$ws.Range("K126").Value = 32664
$ws.Range("L126").Value = 11598.9999
$ws.Range("M126").Value = -30194
$ws.Range("N126").Value = -16538.9999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1350.7391
$ws.Range("I2").Value = 63.42857
$ws.Range("K2").Value = 380.57142
$ws.Range("M2").Value = -267.57142
$ws.Range("H4").Value = 62268410
$ws.Range("I4").Value = 82073736
$ws.Range("K4").Value = 246221208
$ws.Range("M4").Value = -246221096

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 817.2308
$ws.Range("I16").Value = 817.2308
$ws.Range("K16").Value = 817.2308
$ws.Range("M16").Value = -647.2308
$ws.Range("H46").Value = 2282.6667
$ws.Range("I46").Value = 2299.4
$ws.Range("K46").Value = 2299.4
$ws.Range("M46").Value = -2111.4
$ws.Range("H55").Value = 461.1154
$ws.Range("I55").Value = 208.75
$ws.Range("K55").Value = 208.75
$ws.Range("M55").Value = -35.75
$ws.Range("H68").Value = 2005.4688
$ws.Range("I68").Value = 1985.3448
$ws.Range("K68").Value = 1985.3448
$ws.Range("M68").Value = -1236.3448
$ws.Range("H71").Value = 2005.4688
$ws.Range("I71").Value = 1985.3448
$ws.Range("K71").Value = 9926.724
$ws.Range("M71").Value = -6182.724
$ws.Range("H131").Value = 98996.5
$ws.Range("J131").Value = 98996.5
$ws.Range("L131").Value = 98996.5
$ws.Range("N131").Value = -109076.5
$ws.Range("H132").Value = 4212.972
$ws.Range("I132").Value = 2924.739
$ws.Range("J132").Value = 6492.154
$ws.Range("K132").Value = 8774.217000000001
$ws.Range("L132").Value = 19476.462
$ws.Range("M132").Value = -6244.217000000001
$ws.Range("N132").Value = -24536.462

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
$ws.Range("H124").Value = 55000
$ws.Range("J124").Value = 55000
$ws.Range("L124").Value = 55000
$ws.Range("N124").Value = -64820
$ws.Range("H126").Value = 2789.4348
$ws.Range("I126").Value = 2841.625
$ws.Range("J126").Value = 2670.1428
$ws.Range("K126").Value = 8524.875
$ws.Range("L126").Value = 8010.428400000001
$ws.Range("M126").Value = -6054.875
$ws.Range("N126").Value = -12950.4284
